# Adds two new species-observation rows (9 and 10) to the Artfynd sheet,
# matching the data appended in the source export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: Spillkråka (Dryocopus martius) ------------------------------
$ws.Range("A9").Value = 131110701
$ws.Range("B9").Value = 57881
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 100049
$ws.Range("F9").Value = "Spillkråka"
$ws.Range("G9").Value = "Dryocopus martius"
$ws.Range("H9").Value = "(Linnaeus, 1758)"

# "Antal" is exported as text in this sheet, not a number - force text so a
# bare "1" isn't reinterpreted as numeric.
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "1"

$ws.Range("M9").Value = "lockläte, övriga läten"
$ws.Range("P9").Value = "Myckeläng, Mpd"
$ws.Range("Q9").Value = 631192
$ws.Range("R9").Value = 6928048
$ws.Range("S9").Value = 50
$ws.Range("T9").Value = "Västernorrland"
$ws.Range("U9").Value = "Timrå"
$ws.Range("V9").Value = "Medelpad"
$ws.Range("W9").Value = "Tynderö"

# Dates are stored as plain text (yyyy-mm-dd) rather than real date
# serials - force text so Excel's autodetection doesn't convert them.
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2026-02-10"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2026-02-10"

$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = "Mattias Edman"
$ws.Range("AX9").Value = "Mattias Edman"

# --- Row 10: Skinnlav (Leptogium saturninum) -----------------------------
$ws.Range("A10").Value = 131110722
$ws.Range("B10").Value = 80252
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 6456
$ws.Range("F10").Value = "Skinnlav"
$ws.Range("G10").Value = "Leptogium saturninum"
$ws.Range("H10").Value = "(Dicks.) Nyl."

$ws.Range("P10").Value = "Myckeläng, Mpd"
$ws.Range("Q10").Value = 631117
$ws.Range("R10").Value = 6928040
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Västernorrland"
$ws.Range("U10").Value = "Timrå"
$ws.Range("V10").Value = "Medelpad"
$ws.Range("W10").Value = "Tynderö"

$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2026-02-10"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2026-02-10"

$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = "Mattias Edman"
$ws.Range("AX10").Value = "Mattias Edman"

Write-Output "Rows 9 and 10 populated"
